$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as literal text so numeric-looking
# strings like "81.90" or "0.00001057" are not coerced to numbers, matching
# the source data which stores prices as text.
$priceRange = $ws.Range("D2:D50")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.320.85'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.842.41'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '239.78'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '0.6267'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.9986'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.07473'
$ws.Range('E8').Value = '  -1.46%  '
$ws.Range('D9').Value = '0.2894'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '24.43'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '0.07731'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.842.44'
$ws.Range('E12').Value = '  -2.42%  '
$ws.Range('D13').Value = '4.978'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').Value = '0.6791'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '0.00001057'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '81.90'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '6.168'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '29.335.86'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '229.05'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').Value = '12.30'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '0.9988'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '7.498'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '0.9994'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '158.41'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '8.418'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').Value = '0.06586'
$ws.Range('E28').Value = '  +17.52%  '
$ws.Range('D29').Value = '1.415'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').Value = '1.482'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('D31').Value = '4.105'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').Value = '4.086'
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').Value = '1.823'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').Value = '1.139'
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('D35').Value = '0.6949'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '2.578'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '1.264.01'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('D38').Value = '2.830'
$ws.Range('E38').Value = '  +4.07%  '
$ws.Range('D39').Value = '0.01835'
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('D40').Value = '6.792'
$ws.Range('E40').Value = '  +6.92%  '
$ws.Range('D41').Value = '0.9177'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('D42').Value = '0.9978'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '2.003.92'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '66.02'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('E46').Value = '  +2.75%  '
$ws.Range('D47').Value = '7.061'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').Value = '0.1161'
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('D49').Value = '8.959'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = '0.3945'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('E51').Value = '  -0.03%  '

# Restore the default cell style so no stray number-format style lingers
# on these cells (matches original workbook which left them unstyled).
$priceRange.Style = "Normal"
